# Update the "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets to reflect the latest scrape.
#
#   F3: 97 -> 98
#   F4:  1 ->  3
#   F5: 59 -> 60
#   F6:  1 ->  3

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 98
    $ws.Range("F4").Value = 3
    $ws.Range("F5").Value = 60
    $ws.Range("F6").Value = 3
}
